# Weekly update: insert a new data row at row 15 (pushing the existing
# rows 15-28 down to 16-29) and populate it with the new week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 15; Excel shifts rows 15:28 down
# to 16:29 and the sheet's used range grows to A1:T29 automatically.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new record's data.
$ws.Range("A15").Value = 1
$ws.Range("B15").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C15").Value = 'Arica y Parinacota'
$ws.Range("D15").Value = 44540
$ws.Range("E15").Value = 15
$ws.Range("F15").Value = 'Fruta'
$ws.Range("G15").Value = 100103
$ws.Range("H15").Value = 'Frutos de hueso (carozo)'
$ws.Range("I15").Value = 100103004
$ws.Range("J15").Value = 'Durazno'
$ws.Range("K15").Value = 'Polar King'
$ws.Range("L15").Value = 'Segunda'
$ws.Range("M15").Value = 270
$ws.Range("N15").Value = 20000
$ws.Range("O15").Value = 21000
$ws.Range("P15").Value = 20500
$ws.Range("Q15").Value = '$/bandeja 18 kilos granel'
$ws.Range("R15").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S15").Value = 1139
$ws.Range("T15").Value = 18
